# Région 02 - nouvelle version du fichier supergroup.xlsx
# Met a jour les groupes (colonne D), les codes (colonne C) des lignes 2 a 13,
# puis deplace la selection active sur E11 (cellule active lors de
# l'enregistrement du classeur par l'auteur).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repositionnement / redimensionnement de la fenetre du classeur
# (xWindow="270"->"390", windowWidth="27495"->"23655").
$win = $wb.Windows.Item(1)
$win.Left = 390
$win.Width = 23655

# Nouvel ordre / nouveaux libelles pour la colonne D (IdxSG -> DesSG)
# et nouveaux codes pour la colonne C (IdxPar), lignes 2 a 13.
$rows = @(
    @{ Row = 2;  Code = 305020100; Label = "Groupes morgiens d'origine" },
    @{ Row = 3;  Code = 305020200; Label = "Enfance" },
    @{ Row = 4;  Code = 305020300; Label = "Jeunesse" },
    @{ Row = 5;  Code = 305020400; Label = "Aînés" },
    @{ Row = 6;  Code = 305020500; Label = "Terre Nouvelle" },
    @{ Row = 7;  Code = 305020600; Label = "REGISTRES PAROISSIAUX" },
    @{ Row = 8;  Code = 305020700; Label = "Adultes" },
    @{ Row = 9;  Code = 305020800; Label = "St-Nicolas" },
    @{ Row = 10; Code = 305020900; Label = "Culte" },
    @{ Row = 11; Code = 305021000; Label = "Annuaire / responsables-animateurs" },
    @{ Row = 12; Code = 305021100; Label = "Noël" },
    @{ Row = 13; Code = 305021200; Label = "Secrétariat" }
)

foreach ($item in $rows) {
    $r = $item.Row
    # Colonne B (IdxPar racine) reste 2040000000 sur toutes les lignes.
    $ws.Cells.Item($r, 2).Value = 2040000000
    # Colonne C (IdxSG) recoit les nouveaux codes.
    $ws.Cells.Item($r, 3).Value = $item.Code
    # Colonne D (DesSG) recoit les libelles dans le nouvel ordre.
    $ws.Cells.Item($r, 4).Value = $item.Label
}

# La cellule active / selection lors du dernier enregistrement est E11.
$ws.Range("E11").Select()
